# adding averages and more checks
$wb = $excel.ActiveWorkbook

$wsTraining = $wb.Worksheets.Item("Training Dashboard")
$wsExam     = $wb.Worksheets.Item("Exam Dashboard")

# ---------------------------------------------------------------------------
# 1) Header formatting: the big bold title (row 1) and the bold header row
#    (row 2) both switch to a bold WHITE font (the title also drops from
#    14pt down to the normal 11pt) on both sheets.
# ---------------------------------------------------------------------------
foreach ($ws in @($wsTraining, $wsExam)) {
    $lastCol = $ws.UsedRange.Columns.Count

    $ws.Range("A1").Font.Size  = 11
    $ws.Range("A1").Font.Color = 16777215   # white (RGB 255,255,255)

    $headerRange = $ws.Range($ws.Cells.Item(2, 1), $ws.Cells.Item(2, $lastCol))
    $headerRange.Font.Color = 16777215      # white (RGB 255,255,255)
}

# ---------------------------------------------------------------------------
# 2) Training Dashboard: "LAST UPDATE" (col I) moved forward 8 days, from
#    08-Sep-2025 to 16-Sep-2025, so "PERIOD TO EXPIRE" (col H) drops by 8
#    for every data row (3-19).
# ---------------------------------------------------------------------------
for ($row = 3; $row -le 19; $row++) {
    $periodCell = $wsTraining.Cells.Item($row, 8)   # column H
    $periodCell.Value2 = $periodCell.Value2 - 8

    # Leading apostrophe forces literal text so it doesn't get auto-parsed
    # into a date serial (matches the existing "08-Sep-2025" text entries).
    $wsTraining.Cells.Item($row, 9).Value2 = "'16-Sep-2025"   # column I
}

# ---------------------------------------------------------------------------
# 3) Exam Dashboard: widen the COMMENTS column and reword the per-row
#    comment from "OK" to "date is valid" for the seven exam rows (3-9).
# ---------------------------------------------------------------------------
$wsExam.Range("E1").EntireColumn.ColumnWidth = 14.17

for ($row = 3; $row -le 9; $row++) {
    $wsExam.Cells.Item($row, 5).Value2 = "date is valid"
}
